$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 62.41592166666667
$ws.Cells.Item(2, 8).Value = 187.247765
$ws.Cells.Item(2, 9).Value = 0.1654944774607044
$ws.Cells.Item(2, 10).Value = 0.1654944774607044
$ws.Cells.Item(2, 13).Value = 3.795192333333334
$ws.Cells.Item(2, 14).Value = 11.385577
$ws.Cells.Item(2, 15).Value = 0.01044213755712683
$ws.Cells.Item(2, 16).Value = 0.01044213755712683
$ws.Cells.Item(2, 17).Value = 236.8804273872673
$ws.Cells.Item(2, 18).Value = 2131.923846485406
$ws.Cells.Item(2, 19).Value = 0.001728116098589502
$ws.Cells.Item(2, 20).Value = 0.001728116098589502
$ws.Cells.Item(3, 7).Value = 62.41592166666667
$ws.Cells.Item(3, 8).Value = 187.247765
$ws.Cells.Item(3, 9).Value = 0.1654944774607044
$ws.Cells.Item(3, 10).Value = 0.1654944774607044
$ws.Cells.Item(3, 13).Value = 243.3763986666667
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.6696287328350964
$ws.Cells.Item(3, 16).Value = 0.6696287328350964
$ws.Cells.Item(3, 17).Value = 15190.56223469411
$ws.Cells.Item(3, 18).Value = 136715.060112247
$ws.Cells.Item(3, 19).Value = 0.1108198572332179
$ws.Cells.Item(3, 20).Value = 0.1108198572332179
$ws.Cells.Item(4, 7).Value = 62.41592166666667
$ws.Cells.Item(4, 8).Value = 187.247765
$ws.Cells.Item(4, 9).Value = 0.1654944774607044
$ws.Cells.Item(4, 10).Value = 0.1654944774607044
$ws.Cells.Item(4, 13).Value = 29.801371
$ws.Cells.Item(4, 14).Value = 89.404113
$ws.Cells.Item(4, 15).Value = 0.08199584844219236
$ws.Cells.Item(4, 16).Value = 0.08199584844219235
$ws.Cells.Item(4, 17).Value = 1860.080037895272
$ws.Cells.Item(4, 18).Value = 16740.72034105744
$ws.Cells.Item(4, 19).Value = 0.01356986009188774
$ws.Cells.Item(4, 20).Value = 0.01356986009188774
$ws.Cells.Item(5, 7).Value = 62.41592166666667
$ws.Cells.Item(5, 8).Value = 187.247765
$ws.Cells.Item(5, 9).Value = 0.1654944774607044
$ws.Cells.Item(5, 10).Value = 0.1654944774607044
$ws.Cells.Item(5, 13).Value = 86.47679266666667
$ws.Cells.Item(5, 14).Value = 259.430378
$ws.Cells.Item(5, 15).Value = 0.2379332811655844
$ws.Cells.Item(5, 16).Value = 0.2379332811655844
$ws.Cells.Item(5, 17).Value = 5397.528717067242
$ws.Cells.Item(5, 18).Value = 48577.75845360518
$ws.Cells.Item(5, 19).Value = 0.03937664403700926
$ws.Cells.Item(5, 20).Value = 0.03937664403700926
$ws.Cells.Item(6, 7).Value = 164.7897643333334
$ws.Cells.Item(6, 8).Value = 494.369293
$ws.Cells.Item(6, 9).Value = 0.4369365253446571
$ws.Cells.Item(6, 10).Value = 0.436936525344657
$ws.Cells.Item(6, 13).Value = 3.795192333333334
$ws.Cells.Item(6, 14).Value = 11.385577
$ws.Cells.Item(6, 15).Value = 0.01044213755712683
$ws.Cells.Item(6, 16).Value = 0.01044213755712683
$ws.Cells.Item(6, 17).Value = 625.4088502096736
$ws.Cells.Item(6, 18).Value = 5628.679651887062
$ws.Cells.Item(6, 19).Value = 0.004562551301381944
$ws.Cells.Item(6, 20).Value = 0.004562551301381943
$ws.Cells.Item(7, 7).Value = 164.7897643333334
$ws.Cells.Item(7, 8).Value = 494.369293
$ws.Cells.Item(7, 9).Value = 0.4369365253446571
$ws.Cells.Item(7, 10).Value = 0.436936525344657
$ws.Cells.Item(7, 13).Value = 243.3763986666667
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.6696287328350964
$ws.Cells.Item(7, 16).Value = 0.6696287328350964
$ws.Cells.Item(7, 17).Value = 40105.93938057539
$ws.Cells.Item(7, 18).Value = 360953.4544251785
$ws.Cells.Item(7, 19).Value = 0.2925852517959127
$ws.Cells.Item(7, 20).Value = 0.2925852517959127
$ws.Cells.Item(8, 7).Value = 164.7897643333334
$ws.Cells.Item(8, 8).Value = 494.369293
$ws.Cells.Item(8, 9).Value = 0.4369365253446571
$ws.Cells.Item(8, 10).Value = 0.436936525344657
$ws.Cells.Item(8, 13).Value = 29.801371
$ws.Cells.Item(8, 14).Value = 89.404113
$ws.Cells.Item(8, 15).Value = 0.08199584844219236
$ws.Cells.Item(8, 16).Value = 0.08199584844219235
$ws.Cells.Item(8, 17).Value = 4910.960903900234
$ws.Cells.Item(8, 18).Value = 44198.64813510211
$ws.Cells.Item(8, 19).Value = 0.03582698111101865
$ws.Cells.Item(8, 20).Value = 0.03582698111101863
$ws.Cells.Item(9, 7).Value = 164.7897643333334
$ws.Cells.Item(9, 8).Value = 494.369293
$ws.Cells.Item(9, 9).Value = 0.4369365253446571
$ws.Cells.Item(9, 10).Value = 0.436936525344657
$ws.Cells.Item(9, 13).Value = 86.47679266666667
$ws.Cells.Item(9, 14).Value = 259.430378
$ws.Cells.Item(9, 15).Value = 0.2379332811655844
$ws.Cells.Item(9, 16).Value = 0.2379332811655844
$ws.Cells.Item(9, 17).Value = 14250.49028384253
$ws.Cells.Item(9, 18).Value = 128254.4125545828
$ws.Cells.Item(9, 19).Value = 0.1039617411363438
$ws.Cells.Item(9, 20).Value = 0.1039617411363438
$ws.Cells.Item(10, 7).Value = 57.486235
$ws.Cells.Item(10, 8).Value = 172.458705
$ws.Cells.Item(10, 9).Value = 0.1524235190071549
$ws.Cells.Item(10, 10).Value = 0.1524235190071549
$ws.Cells.Item(10, 13).Value = 3.795192333333334
$ws.Cells.Item(10, 14).Value = 11.385577
$ws.Cells.Item(10, 15).Value = 0.01044213755712683
$ws.Cells.Item(10, 16).Value = 0.01044213755712683
$ws.Cells.Item(10, 17).Value = 218.1713183441984
$ws.Cells.Item(10, 18).Value = 1963.541865097785
$ws.Cells.Item(10, 19).Value = 0.001591627352414048
$ws.Cells.Item(10, 20).Value = 0.001591627352414048
$ws.Cells.Item(11, 7).Value = 57.486235
$ws.Cells.Item(11, 8).Value = 172.458705
$ws.Cells.Item(11, 9).Value = 0.1524235190071549
$ws.Cells.Item(11, 10).Value = 0.1524235190071549
$ws.Cells.Item(11, 13).Value = 243.3763986666667
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.6696287328350964
$ws.Cells.Item(11, 16).Value = 0.6696287328350964
$ws.Cells.Item(11, 17).Value = 13990.79284720569
$ws.Cells.Item(11, 18).Value = 125917.1356248512
$ws.Cells.Item(11, 19).Value = 0.1020671678870274
$ws.Cells.Item(11, 20).Value = 0.1020671678870274
$ws.Cells.Item(12, 7).Value = 57.486235
$ws.Cells.Item(12, 8).Value = 172.458705
$ws.Cells.Item(12, 9).Value = 0.1524235190071549
$ws.Cells.Item(12, 10).Value = 0.1524235190071549
$ws.Cells.Item(12, 13).Value = 29.801371
$ws.Cells.Item(12, 14).Value = 89.404113
$ws.Cells.Item(12, 15).Value = 0.08199584844219236
$ws.Cells.Item(12, 16).Value = 0.08199584844219235
$ws.Cells.Item(12, 17).Value = 1713.168616628185
$ws.Cells.Item(12, 18).Value = 15418.51754965366
$ws.Cells.Item(12, 19).Value = 0.0124980957635363
$ws.Cells.Item(12, 20).Value = 0.0124980957635363
$ws.Cells.Item(13, 7).Value = 57.486235
$ws.Cells.Item(13, 8).Value = 172.458705
$ws.Cells.Item(13, 9).Value = 0.1524235190071549
$ws.Cells.Item(13, 10).Value = 0.1524235190071549
$ws.Cells.Item(13, 13).Value = 86.47679266666667
$ws.Cells.Item(13, 14).Value = 259.430378
$ws.Cells.Item(13, 15).Value = 0.2379332811655844
$ws.Cells.Item(13, 16).Value = 0.2379332811655844
$ws.Cells.Item(13, 17).Value = 4971.225225282276
$ws.Cells.Item(13, 18).Value = 44741.0270275405
$ws.Cells.Item(13, 19).Value = 0.03626662800417718
$ws.Cells.Item(13, 20).Value = 0.03626662800417718
$ws.Cells.Item(14, 7).Value = 92.45614233333333
$ws.Cells.Item(14, 8).Value = 277.368427
$ws.Cells.Item(14, 9).Value = 0.2451454781874835
$ws.Cells.Item(14, 10).Value = 0.2451454781874835
$ws.Cells.Item(14, 13).Value = 3.795192333333334
$ws.Cells.Item(14, 14).Value = 11.385577
$ws.Cells.Item(14, 15).Value = 0.01044213755712683
$ws.Cells.Item(14, 16).Value = 0.01044213755712683
$ws.Cells.Item(14, 17).Value = 350.8888425530421
$ws.Cells.Item(14, 18).Value = 3157.999582977379
$ws.Cells.Item(14, 19).Value = 0.002559842804741339
$ws.Cells.Item(14, 20).Value = 0.002559842804741339
$ws.Cells.Item(15, 7).Value = 92.45614233333333
$ws.Cells.Item(15, 8).Value = 277.368427
$ws.Cells.Item(15, 9).Value = 0.2451454781874835
$ws.Cells.Item(15, 10).Value = 0.2451454781874835
$ws.Cells.Item(15, 13).Value = 243.3763986666667
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.6696287328350964
$ws.Cells.Item(15, 16).Value = 0.6696287328350964
$ws.Cells.Item(15, 17).Value = 22501.64295569941
$ws.Cells.Item(15, 18).Value = 202514.7866012947
$ws.Cells.Item(15, 19).Value = 0.1641564559189383
$ws.Cells.Item(15, 20).Value = 0.1641564559189383
$ws.Cells.Item(16, 7).Value = 92.45614233333333
$ws.Cells.Item(16, 8).Value = 277.368427
$ws.Cells.Item(16, 9).Value = 0.2451454781874835
$ws.Cells.Item(16, 10).Value = 0.2451454781874835
$ws.Cells.Item(16, 13).Value = 29.801371
$ws.Cells.Item(16, 14).Value = 89.404113
$ws.Cells.Item(16, 15).Value = 0.08199584844219236
$ws.Cells.Item(16, 16).Value = 0.08199584844219235
$ws.Cells.Item(16, 17).Value = 2755.319798904472
$ws.Cells.Item(16, 18).Value = 24797.87819014025
$ws.Cells.Item(16, 19).Value = 0.02010091147574967
$ws.Cells.Item(16, 20).Value = 0.02010091147574967
$ws.Cells.Item(17, 7).Value = 92.45614233333333
$ws.Cells.Item(17, 8).Value = 277.368427
$ws.Cells.Item(17, 9).Value = 0.2451454781874835
$ws.Cells.Item(17, 10).Value = 0.2451454781874835
$ws.Cells.Item(17, 13).Value = 86.47679266666667
$ws.Cells.Item(17, 14).Value = 259.430378
$ws.Cells.Item(17, 15).Value = 0.2379332811655844
$ws.Cells.Item(17, 16).Value = 0.2379332811655844
$ws.Cells.Item(17, 17).Value = 7995.31065131949
$ws.Cells.Item(17, 18).Value = 71957.79586187541
$ws.Cells.Item(17, 19).Value = 0.05832826798805415
$ws.Cells.Item(17, 20).Value = 0.05832826798805415
